$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest scrape of the cryptos list: updated Price / Volume(1h) cells,
# plus a reordering of the Maker/Stacks rows (39 <-> 40).
$updates = @(
    @{ Cell = "D2"; Value = "58.784.27" },
    @{ Cell = "E2"; Value = "  +1.41%  " },
    @{ Cell = "D3"; Value = "3.157.74" },
    @{ Cell = "E3"; Value = "  +1.08%  " },
    @{ Cell = "E4"; Value = "  +0.01%  " },
    @{ Cell = "D5"; Value = "532.47" },
    @{ Cell = "E5"; Value = "  +0.48%  " },
    @{ Cell = "D6"; Value = "140.00" },
    @{ Cell = "E6"; Value = "  +1.04%  " },
    @{ Cell = "E7"; Value = "  +0.10%  " },
    @{ Cell = "D8"; Value = "0.530" },
    @{ Cell = "E8"; Value = "  +15.23%  " },
    @{ Cell = "D9"; Value = "7.32" },
    @{ Cell = "E9"; Value = "  +0.50%  " },
    @{ Cell = "D10"; Value = "0.429" },
    @{ Cell = "E10"; Value = "  +5.53%  " },
    @{ Cell = "E11"; Value = "  +2.90%  " },
    @{ Cell = "E12"; Value = "  +2.74%  " },
    @{ Cell = "D13"; Value = "3.698.09" },
    @{ Cell = "E13"; Value = "  +1.16%  " },
    @{ Cell = "D14"; Value = "25.89" },
    @{ Cell = "E14"; Value = "  +1.62%  " },
    @{ Cell = "E15"; Value = "  +4.99%  " },
    @{ Cell = "D16"; Value = "58.808.99" },
    @{ Cell = "E16"; Value = "  +1.44%  " },
    @{ Cell = "D17"; Value = "6.23" },
    @{ Cell = "E17"; Value = "  +4.15%  " },
    @{ Cell = "D18"; Value = "3.152.30" },
    @{ Cell = "E18"; Value = "  +1.02%  " },
    @{ Cell = "D19"; Value = "13.02" },
    @{ Cell = "E19"; Value = "  +2.92%  " },
    @{ Cell = "D20"; Value = "8.22" },
    @{ Cell = "E20"; Value = "  +1.54%  " },
    @{ Cell = "D21"; Value = "370.67" },
    @{ Cell = "E21"; Value = "  +5.36%  " },
    @{ Cell = "D22"; Value = "5.80" },
    @{ Cell = "E22"; Value = "  +1.82%  " },
    @{ Cell = "D23"; Value = "0.999" },
    @{ Cell = "E23"; Value = "  -0.04%  " },
    @{ Cell = "D24"; Value = "69.99" },
    @{ Cell = "E24"; Value = "  +1.54%  " },
    @{ Cell = "E25"; Value = "  +3.05%  " },
    @{ Cell = "E26"; Value = "  +1.13%  " },
    @{ Cell = "D27"; Value = "1.00" },
    @{ Cell = "E27"; Value = "  -0.03%  " },
    @{ Cell = "E28"; Value = "  +13.35%  " },
    @{ Cell = "E29"; Value = "  -1.83%  " },
    @{ Cell = "E30"; Value = "  +1.29%  " },
    @{ Cell = "D31"; Value = "6.12" },
    @{ Cell = "E31"; Value = "  +0.48%  " },
    @{ Cell = "D32"; Value = "21.99" },
    @{ Cell = "E32"; Value = "  +3.39%  " },
    @{ Cell = "D33"; Value = "5.21" },
    @{ Cell = "E33"; Value = "  +5.24%  " },
    @{ Cell = "E34"; Value = "  +1.76%  " },
    @{ Cell = "D35"; Value = "159.27" },
    @{ Cell = "E35"; Value = "  +0.26%  " },
    @{ Cell = "D36"; Value = "6.29" },
    @{ Cell = "E36"; Value = "  +4.09%  " },
    @{ Cell = "D37"; Value = "1.35" },
    @{ Cell = "E37"; Value = "  +7.48%  " },
    @{ Cell = "D38"; Value = "25.35" },
    @{ Cell = "E38"; Value = "  -2.92%  " },
    @{ Cell = "B39"; Value = "Stacks" },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx" },
    @{ Cell = "D39"; Value = "1.67" },
    @{ Cell = "E39"; Value = "  +0.63%  " },
    @{ Cell = "B40"; Value = "Maker" },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" },
    @{ Cell = "D40"; Value = "2.645.70" },
    @{ Cell = "E40"; Value = "  +10.33%  " },
    @{ Cell = "D41"; Value = "0.0682" },
    @{ Cell = "E41"; Value = "  +2.16%  " },
    @{ Cell = "E42"; Value = "  +6.11%  " },
    @{ Cell = "D43"; Value = "38.78" },
    @{ Cell = "E43"; Value = "  +3.42%  " },
    @{ Cell = "D44"; Value = "0.711" },
    @{ Cell = "E44"; Value = "  +2.00%  " },
    @{ Cell = "E45"; Value = "  +7.54%  " },
    @{ Cell = "D46"; Value = "0.999" },
    @{ Cell = "E46"; Value = "  +0.02%  " },
    @{ Cell = "D47"; Value = "3.196.83" },
    @{ Cell = "E47"; Value = "  +1.09%  " },
    @{ Cell = "E48"; Value = "  +13.83%  " },
    @{ Cell = "D49"; Value = "0.986" },
    @{ Cell = "E49"; Value = "  +1.85%  " },
    @{ Cell = "D50"; Value = "6.20" },
    @{ Cell = "E50"; Value = "  +2.98%  " },
    @{ Cell = "D51"; Value = "20.34" },
    @{ Cell = "E51"; Value = "  +3.25%  " }
)

foreach ($u in $updates) {
    # Force text so numeric-looking strings (e.g. "532.47", "0.530")
    # keep their exact original formatting instead of becoming numbers.
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
    # Restore default styling so no stray number-format style sticks.
    $ws.Range($u.Cell).Style = "Normal"
}
